$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "acronym" header in E1
$ws.Cells.Item(1, 5).Value = "acronym"

# Fill in the acronym values for each data row (E2:E11)
$acronyms = @("area_a", "area_b", "area_c", "area_a", "area_b", "area_c", "area_a", "area_b", "area_c", "area_a")
for ($i = 0; $i -lt $acronyms.Length; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $acronyms[$i]
}

# New column's data cells get a wrap-text style
$ws.Range("E2:E11").WrapText = $true

# Leave the new column selected, matching the state left after the import
$ws.Range("E1:E11").Select()
